# Applies the numeric corrections to the Leve profit tables (H:N columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# price-refresh run.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 49386
$ws.Range("I26").Value = 1900
$ws.Range("J26").Value = 61257.5
$ws.Range("K26").Value = 1900
$ws.Range("L26").Value = 61257.5
$ws.Range("M26").Value = -1556
$ws.Range("N26").Value = -61945.5
$ws.Range("H38").Value = 62.125
$ws.Range("I38").Value = 32.933334
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 98.800002
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 273.199998
$ws.Range("N38").Value = -2244
$ws.Range("H64").Value = 3461.25
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3527.1428
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3527.1428
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4023.1428
$ws.Range("H67").Value = 3461.25
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3527.1428
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3527.1428
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5243.1428
$ws.Range("H112").Value = 1474.4667
$ws.Range("I112").Value = 1100.5
$ws.Range("J112").Value = 1532
$ws.Range("K112").Value = 3301.5
$ws.Range("L112").Value = 4596
$ws.Range("M112").Value = -2193.5
$ws.Range("N112").Value = -6812
$ws.Range("H137").Value = 1729.25
$ws.Range("I137").Value = 1126.2
$ws.Range("J137").Value = 2734.3333
$ws.Range("K137").Value = 3378.6
$ws.Range("L137").Value = 8202.999899999999
$ws.Range("M137").Value = -828.6000000000004
$ws.Range("N137").Value = -13302.9999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 4413.875
$ws.Range("I26").Value = 3551.8333
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 3551.8333
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = -3221.8333
$ws.Range("N26").Value = -7660
$ws.Range("H32").Value = 18259.521
$ws.Range("I32").Value = 22005.352
$ws.Range("J32").Value = 7613.4736
$ws.Range("K32").Value = 22005.352
$ws.Range("L32").Value = 7613.4736
$ws.Range("M32").Value = -21718.352
$ws.Range("N32").Value = -8187.4736
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2480
$ws.Range("N39").ClearContents()
$ws.Range("H132").Value = 5629.5
$ws.Range("I132").Value = 6684.619
$ws.Range("J132").Value = 3167.5557
$ws.Range("K132").Value = 20053.857
$ws.Range("L132").Value = 9502.667099999999
$ws.Range("M132").Value = -17523.857
$ws.Range("N132").Value = -14562.6671

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 251.5
$ws.Range("I7").Value = 135.33333
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 135.33333
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -22.33332999999999
$ws.Range("N7").Value = -826
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 2000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 2000
$ws.Range("N56").Value = -3478
$ws.Range("H64").Value = 370.54544
$ws.Range("I64").Value = 355.33334
$ws.Range("J64").Value = 376.25
$ws.Range("K64").Value = 355.33334
$ws.Range("L64").Value = 376.25
$ws.Range("M64").Value = -130.33334
$ws.Range("N64").Value = -826.25
$ws.Range("H67").Value = 370.54544
$ws.Range("I67").Value = 355.33334
$ws.Range("J67").Value = 376.25
$ws.Range("K67").Value = 355.33334
$ws.Range("L67").Value = 376.25
$ws.Range("M67").Value = 424.66666
$ws.Range("N67").Value = -1936.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 350
$ws.Range("I4").Value = 350
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 350
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -238
$ws.Range("H5").Value = 234.09091
$ws.Range("I5").Value = 209.57143
$ws.Range("J5").Value = 277
$ws.Range("K5").Value = 209.57143
$ws.Range("L5").Value = 277
$ws.Range("M5").Value = -97.57142999999999
$ws.Range("N5").Value = -501
$ws.Range("H12").Value = 1072.6666
$ws.Range("I12").Value = 1259
$ws.Range("J12").Value = 700
$ws.Range("K12").Value = 1259
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = -1089
$ws.Range("N12").Value = -1040
$ws.Range("H25").Value = 1900
$ws.Range("I25").Value = 1900
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1900
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1726
$ws.Range("H35").Value = 37048
$ws.Range("I35").Value = 833.3333
$ws.Range("J35").Value = 50628.5
$ws.Range("K35").Value = 833.3333
$ws.Range("L35").Value = 50628.5
$ws.Range("M35").Value = -539.3333
$ws.Range("N35").Value = -51216.5
$ws.Range("H59").Value = 36083.332
$ws.Range("I59").Value = 20250
$ws.Range("J59").Value = 44000
$ws.Range("K59").Value = 20250
$ws.Range("L59").Value = 44000
$ws.Range("M59").Value = -19105
$ws.Range("N59").Value = -46290

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 284.33334
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 257.375
$ws.Range("K15").Value = 1500
$ws.Range("L15").Value = 772.125
$ws.Range("M15").Value = -1360
$ws.Range("N15").Value = -1052.125
$ws.Range("H17").Value = 236
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 172
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 516
$ws.Range("M17").Value = -731
$ws.Range("N17").Value = -854
$ws.Range("H110").Value = 25633.334
$ws.Range("I110").Value = 3800
$ws.Range("J110").Value = 30000
$ws.Range("K110").Value = 11400
$ws.Range("L110").Value = 90000
$ws.Range("M110").Value = -7310
$ws.Range("N110").Value = -98180
$ws.Range("H122").Value = 1545.7894
$ws.Range("I122").Value = 1774.125
$ws.Range("J122").Value = 1379.7273
$ws.Range("K122").Value = 15967.125
$ws.Range("L122").Value = 12417.5457
$ws.Range("M122").Value = -13517.125
$ws.Range("N122").Value = -17317.5457
$ws.Range("H131").Value = 3271922.5
$ws.Range("I131").Value = 10341.363
$ws.Range("J131").Value = 4831809.5
$ws.Range("K131").Value = 31024.089
$ws.Range("L131").Value = 14495428.5
$ws.Range("M131").Value = -25984.089
$ws.Range("N131").Value = -14505508.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H102").Value = 1258.4166
$ws.Range("I102").Value = 1189
$ws.Range("J102").Value = 1466.6666
$ws.Range("K102").Value = 1189
$ws.Range("L102").Value = 1466.6666
$ws.Range("M102").Value = 433
$ws.Range("N102").Value = -4710.6666
$ws.Range("H132").Value = 35326.168
$ws.Range("I132").Value = 44977.086
$ws.Range("J132").Value = 3616
$ws.Range("K132").Value = 134931.258
$ws.Range("L132").Value = 10848
$ws.Range("M132").Value = -132401.258
$ws.Range("N132").Value = -15908

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7125.7646
$ws.Range("I16").Value = 799.0769
$ws.Range("J16").Value = 27687.5
$ws.Range("K16").Value = 799.0769
$ws.Range("L16").Value = 27687.5
$ws.Range("M16").Value = -629.0769
$ws.Range("N16").Value = -28027.5
$ws.Range("H122").Value = 4071.7273
$ws.Range("I122").Value = 9242
$ws.Range("J122").Value = 2922.7778
$ws.Range("K122").Value = 27726
$ws.Range("L122").Value = 8768.3334
$ws.Range("M122").Value = -25276
$ws.Range("N122").Value = -13668.3334
$ws.Range("H132").Value = 1598.0588
$ws.Range("I132").Value = 1385.1025
$ws.Range("J132").Value = 2290.1667
$ws.Range("K132").Value = 4155.3075
$ws.Range("L132").Value = 6870.500100000001
$ws.Range("M132").Value = -1625.3075
$ws.Range("N132").Value = -11930.5001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21240
$ws.Range("H132").Value = 1081.8235
$ws.Range("I132").Value = 923.8082
$ws.Range("J132").Value = 2043.0834
$ws.Range("K132").Value = 2771.4246
$ws.Range("L132").Value = 6129.2502
$ws.Range("M132").Value = -241.4246000000003
$ws.Range("N132").Value = -11189.2502
$ws.Range("H136").Value = 10110.533
$ws.Range("I136").Value = 10912.154
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 32736.462
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = -30186.462
$ws.Range("N136").Value = -19800
